{"js": "// Change the final chapter heading from \"CHAPTER 14\" to \"CHAPTER 15\".\n// (commit message: \"chapter 13 project and chapter 15 demo\")\n//\n// The heading lives in its own paragraph, styled bold + single underline,\n// and is the ONLY occurrence of the literal text \"CHAPTER 14\" in the body.\nconst results = context.document.body.search(\"CHAPTER 14\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Replacing in place keeps the run's existing formatting (bold + single\n  // underline) exactly as it was - only the visible text changes.\n  results.items[i].insertText(\"CHAPTER 15\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Change the final chapter heading from \"CHAPTER 14\" to \"CHAPTER 15\".\n# (commit message: \"chapter 13 project and chapter 15 demo\")\n#\n# The heading lives in its own paragraph, styled bold + single underline,\n# and is the ONLY occurrence of the literal text \"CHAPTER 14\" in the body.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"CHAPTER 14\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"CHAPTER 15\"\n\n# wdFindContinue = 1 (keep searching the whole story),\n# wdReplaceAll = 2 (replace every match - there is exactly one).\n$rng.Find.Execute(\n    [ref]\"CHAPTER 14\",   # FindText\n    [ref]$true,          # MatchCase\n    [ref]$true,          # MatchWholeWord\n    [ref]$false,         # MatchWildcards\n    [ref]$false,         # MatchSoundsLike\n    [ref]$false,         # MatchAllWordForms\n    [ref]$true,          # Forward\n    [ref]1,              # Wrap (wdFindContinue)\n    [ref]$false,         # Format\n    [ref]\"CHAPTER 15\",   # ReplaceWith\n    [ref]2               # Replace (wdReplaceAll)\n) | Out-Null\n"}
